$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.388.69'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").Value = '1.871.71'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.021'
$ws.Range("E4").Value = '  +1.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.33'
$ws.Range("E5").Value = '  +0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.018'
$ws.Range("E6").Value = '  +0.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5105'
$ws.Range("E7").Value = '  -0.67%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3949'
$ws.Range("E8").Value = '  +1.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08451'
$ws.Range("E9").Value = '  +0.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.108'
$ws.Range("E10").Value = '  -1.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.95'
$ws.Range("E11").Value = '  +0.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.245'
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").Value = '1.867.61'
$ws.Range("E13").Value = '  -1.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.44'
$ws.Range("E14").Value = '  -0.80%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.224'
$ws.Range("E15").Value = '  -0.57%  '
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.020'
$ws.Range("E16").Value = '  +1.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001110'
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.84'
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06763'
$ws.Range("E19").Value = '  +1.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.72'
$ws.Range("E20").Value = '  -0.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.018'
$ws.Range("E21").Value = '  +0.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.938'
$ws.Range("E22").Value = '  -1.52%  '
$ws.Range("D23").Value = '28.438.40'
$ws.Range("E23").Value = '  +0.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.15'
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.291'
$ws.Range("E25").Value = '  +0.54%  '
$ws.Range("D26").Value = '2.081.98'
$ws.Range("E26").Value = '  -0.97%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.70'
$ws.Range("E27").Value = '  +0.91%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.75'
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.352'
$ws.Range("E29").Value = '  -4.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.11'
$ws.Range("E30").Value = '  +1.26%  '
$ws.Range("E31").Value = '  -0.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.039'
$ws.Range("E32").Value = '  -0.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.756'
$ws.Range("E33").Value = '  -2.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.634'
$ws.Range("E34").Value = '  +0.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02436'
$ws.Range("E35").Value = '  -0.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06452'
$ws.Range("E36").Value = '  -1.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2172'
$ws.Range("E37").Value = '  -1.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.824'
$ws.Range("E38").Value = '  -6.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.259'
$ws.Range("E39").Value = '  +1.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.183'
$ws.Range("E40").Value = '  -1.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6360'
$ws.Range("E41").Value = '  -2.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.983'
$ws.Range("E42").Value = '  -0.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.21'
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6015'
$ws.Range("E44").Value = '  -1.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.01'
$ws.Range("E45").Value = '  -1.30%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.703'
$ws.Range("E46").Value = '  +0.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.210'
$ws.Range("E47").Value = '  -5.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.989'
$ws.Range("E48").Value = '  -1.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.207'
$ws.Range("E49").Value = '  -2.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '121.56'
$ws.Range("E50").Value = '  +0.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06842'
$ws.Range("E51").Value = '  -1.03%  '
